$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-15 08:50:01"
$wsZhCn.Range("G5").Value = "2016-02-15 08:50:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-15 08:50:17"
$wsDeDe.Range("G5").Value = "2016-02-15 08:51:18"
